$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated symbol list with latest cryptocurrency prices / volume figures.
# All target cells store plain text (not numbers), so the number format is
# temporarily forced to Text ("@") while assigning the value, then restored
# to General, to stop Excel auto-converting numeric-looking strings.
$cellUpdates = @{
    'D2' = '305.25'
    'E2' = '0.17%'
    'D3' = '35.77'
    'E3' = '-0.31%'
    'D4' = '5.039'
    'E4' = '-0.85%'
    'D5' = '0.08033'
    'E5' = '-0.23%'
    'D6' = '1.864'
    'E6' = '-2.95%'
    'D7' = '4.145'
    'E7' = '-0.76%'
    'D8' = '7.792'
    'E8' = '0.64%'
    'D9' = '0.9197'
    'D10' = '0.1272'
    'E10' = '-7.40%'
    'D11' = '0.1915'
    'E11' = '0.52%'
    'D12' = '0.09090'
    'E12' = '-0.91%'
    'D13' = '0.03452'
    'E13' = '-3.99%'
    'D14' = '0.09857'
    'E14' = '0.33%'
    'D15' = '0.001405'
    'E15' = '-0.62%'
    'D16' = '0.006163'
    'E16' = '4.26%'
    'D17' = '3.820'
    'E17' = '7.39%'
    'E18' = '12.62%'
    'E19' = '-1.01%'
    'B20' = 'MCDex'
    'C20' = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
    'D20' = '5.237'
    'E20' = '6.82%'
    'B21' = 'ProBitToken'
    'C21' = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
    'D21' = '0.1321'
    'E21' = '-0.88%'
    'E22' = '-8.19%'
    'E23' = '-0.68%'
    'D24' = '0.001235'
    'E24' = '1.07%'
    'D25' = '0.004616'
    'E25' = '-3.80%'
    'D27' = '0.0001252'
    'E27' = '-15.60%'
    'E28' = '41.87%'
    'D39' = '0.01945'
    'E39' = '-0.88%'
    'D40' = '0.05285'
    'E40' = '7.52%'
    'D41' = '0.007634'
    'E41' = '-0.07%'
    'E42' = '10.49%'
    'D43' = '0.1354'
    'E43' = '-1.57%'
    'D44' = '0.002173'
    'E44' = '3.26%'
    'D45' = '0.009625'
    'E45' = '-15.13%'
    'D46' = '0.00006103'
    'E46' = '-4.80%'
    'E47' = '-0.08%'
    'D48' = '64.97'
    'E48' = '0.47%'
    'E49' = '39.19%'
    'D50' = '0.00002103'
    'E50' = '-0.08%'
    'D51' = '0.0002003'
    'E51' = '-0.08%'
}

foreach ($addr in $cellUpdates.Keys) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $cellUpdates[$addr]
    $rng.NumberFormat = "General"
}
